$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update formulas (B12 must change before D10 to avoid a transient circular reference,
# matching the dependency resolution order of the real edit).
$ws.Range("B12").Formula = "=C12-C10"
$ws.Range("D10").Formula = "=D12-B12"

# Highlight the newly-updated cells with a yellow fill (keeps their existing thin border).
$cells = @("D11", "E11", "F11", "F12", "C13", "D13")
foreach ($addr in $cells) {
    $ws.Range($addr).Interior.Color = 65535
}

# Move the active selection to H10, matching the saved view state.
$ws.Range("H10").Select()
